$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "0.289"; New = "0.317" },
    @{ Old = "0.188"; New = "0.194" },
    @{ Old = "0.234"; New = "0.243" },
    @{ Old = "0.128"; New = "0.156" },
    @{ Old = "0.232"; New = "0.224" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.New, 2)
}
